$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new columns before column D (shifts old D:K -> F:M)
$ws.Columns("D:E").Insert()

# Copy number formats from column F (the old column D, now shifted) into new D:E
# so the new columns inherit the same date / accounting-number styles.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarter columns (D, E) and refresh the corrected
# historical values that moved into H and I (FY2017 Q4 / Q3 restatements).
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("H7").Value = 43100
$ws.Range("I7").Value = 43008
$ws.Range("D8").Value = 1077300
$ws.Range("E8").Value = 1181400
$ws.Range("H8").Value = 990500
$ws.Range("I8").Value = 1002400
$ws.Range("D9").Value = 798500
$ws.Range("E9").Value = 867100
$ws.Range("H9").Value = 739400
$ws.Range("I9").Value = 717800
$ws.Range("D10").Value = 278800
$ws.Range("E10").Value = 314300
$ws.Range("H10").Value = 251100
$ws.Range("I10").Value = 284600
$ws.Range("D12").Value = 14700
$ws.Range("E12").Value = 12300
$ws.Range("H12").Value = 13800
$ws.Range("I12").Value = 12400
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = "NA"
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("D17").Value = 962700
$ws.Range("E17").Value = 1041400
$ws.Range("H17").Value = 897100
$ws.Range("I17").Value = 867600
$ws.Range("D18").Value = 114600
$ws.Range("E18").Value = 140000
$ws.Range("H18").Value = 93400
$ws.Range("I18").Value = 134800
$ws.Range("D20").Value = 400
$ws.Range("E20").Value = -1800
$ws.Range("H20").Value = 300
$ws.Range("I20").Value = -3300
$ws.Range("D21").Value = 162400
$ws.Range("E21").Value = 184100
$ws.Range("H21").Value = 142100
$ws.Range("I21").Value = 173100
$ws.Range("D22").Value = 15500
$ws.Range("E22").Value = 16300
$ws.Range("H22").Value = 12300
$ws.Range("I22").Value = 7900
$ws.Range("D23").Value = 99500
$ws.Range("E23").Value = 121900
$ws.Range("H23").Value = 81400
$ws.Range("I23").Value = 123600
$ws.Range("D24").Value = 15000
$ws.Range("E24").Value = 23700
$ws.Range("H24").Value = 30200
$ws.Range("I24").Value = 44500
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("D26").Value = 84500
$ws.Range("E26").Value = 98200
$ws.Range("H26").Value = 51200
$ws.Range("I26").Value = 79100
$ws.Range("D27").Value = 86300
$ws.Range("E27").Value = 97400
$ws.Range("H27").Value = 50100
$ws.Range("I27").Value = 78700
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("D29").Value = 3600
$ws.Range("E29").Value = 1500
$ws.Range("H29").Value = 63900
$ws.Range("I29").Value = 7200
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("D32").Value = -400
$ws.Range("E32").Value = 1800
$ws.Range("H32").Value = -300
$ws.Range("I32").Value = 3300
$ws.Range("D33").Value = 89900
$ws.Range("E33").Value = 98900
$ws.Range("H33").Value = 114000
$ws.Range("I33").Value = 85900
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("D35").Value = 89900
$ws.Range("E35").Value = 98900
$ws.Range("H35").Value = 114000
$ws.Range("I35").Value = 85900
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("H38").Value = 43100
$ws.Range("I38").Value = 43008
$ws.Range("D41").Value = 803600
$ws.Range("E41").Value = 780500
$ws.Range("H41").Value = 378300
$ws.Range("I41").Value = 147600
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("D43").Value = 698300
$ws.Range("E43").Value = 850500
$ws.Range("H43").Value = 625700
$ws.Range("I43").Value = 690800
$ws.Range("D44").Value = 457500
$ws.Range("E44").Value = 489100
$ws.Range("H44").Value = 448800
$ws.Range("I44").Value = 454800
$ws.Range("D45").Value = 97300
$ws.Range("E45").Value = 81300
$ws.Range("H45").Value = 191800
$ws.Range("I45").Value = 68500
$ws.Range("D46").Value = 2056700
$ws.Range("E46").Value = 2201400
$ws.Range("H46").Value = 1644600
$ws.Range("I46").Value = 1361700
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("D48").Value = 760100
$ws.Range("E48").Value = 763000
$ws.Range("H48").Value = 1512000
$ws.Range("I48").Value = 698500
$ws.Range("D49").Value = 2409500
$ws.Range("E49").Value = 2436800
$ws.Range("H49").Value = 2517100
$ws.Range("I49").Value = 2231400
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("D52").Value = 22900
$ws.Range("E52").Value = 31200
$ws.Range("H52").Value = 407000
$ws.Range("I52").Value = 26300
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("D54").Value = 5249200
$ws.Range("E54").Value = 5432400
$ws.Range("H54").Value = 5299800
$ws.Range("I54").Value = 4317900
$ws.Range("D57").Value = 312100
$ws.Range("E57").Value = 362700
$ws.Range("H57").Value = 332100
$ws.Range("I57").Value = 328900
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("D59").Value = 283500
$ws.Range("E59").Value = 310000
$ws.Range("H59").Value = 604900
$ws.Range("I59").Value = 306700
$ws.Range("D60").Value = 595600
$ws.Range("E60").Value = 672700
$ws.Range("H60").Value = 658600
$ws.Range("I60").Value = 635600
$ws.Range("D61").Value = 1587800
$ws.Range("E61").Value = 1587400
$ws.Range("H61").Value = 1586200
$ws.Range("I61").Value = 781900
$ws.Range("D62").Value = 468400
$ws.Range("E62").Value = 466600
$ws.Range("H62").Value = 865400
$ws.Range("I62").Value = 464100
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("D66").Value = 2651800
$ws.Range("E66").Value = 2726700
$ws.Range("H66").Value = 2771500
$ws.Range("I66").Value = 1881600
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("D72").Value = 3351400
$ws.Range("E72").Value = 3287100
$ws.Range("H72").Value = 2820800
$ws.Range("I72").Value = 2728800
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("D76").Value = 2597400
$ws.Range("E76").Value = 2705700
$ws.Range("H76").Value = 2528300
$ws.Range("I76").Value = 2436300
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("H80").Value = 43100
$ws.Range("I80").Value = 43008
$ws.Range("D81").Value = 89900
$ws.Range("E81").Value = 98900
$ws.Range("H81").Value = 114000
$ws.Range("I81").Value = 85900
$ws.Range("D83").Value = 47400
$ws.Range("E83").Value = 45900
$ws.Range("H83").Value = 48400
$ws.Range("I83").Value = 41600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("D89").Value = 236600
$ws.Range("E89").Value = 105400
$ws.Range("H89").Value = 159100
$ws.Range("I89").Value = 164900
$ws.Range("D91").Value = -24600
$ws.Range("E91").Value = -29200
$ws.Range("H91").Value = -54100
$ws.Range("I91").Value = -39700
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("D94").Value = -24600
$ws.Range("E94").Value = -20300
$ws.Range("H94").Value = -708700
$ws.Range("I94").Value = -93700
$ws.Range("D96").Value = -23800
$ws.Range("E96").Value = -24100
$ws.Range("H96").Value = -23100
$ws.Range("I96").Value = -23200
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("D100").Value = -188200
$ws.Range("E100").Value = -66600
$ws.Range("H100").Value = 781700
$ws.Range("I100").Value = -64100
$ws.Range("D101").Value = -700
$ws.Range("E101").Value = -400
$ws.Range("H101").Value = -100
$ws.Range("I101").Value = 700
$ws.Range("D102").Value = 23100
$ws.Range("E102").Value = 18100
$ws.Range("H102").Value = 232000
$ws.Range("I102").Value = 7800
